$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:K3").NumberFormat = "@"

$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " September 19 2020"
$ws.Range("C3").Value = "Super Kings won by 5 wickets (with 4 balls remaining)"
$ws.Range("D3").Value = "Mumbai Indians"
$ws.Range("E3").Value = "Chennai Super Kings"
$ws.Range("F3").Value = "Jasprit Bumrah "
$ws.Range("G3").Value = "5"
$ws.Range("H3").Value = "3"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "166.66"

$ws.Range("A3:K3").ClearFormats()
